$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order A..AL (1-based index matches PS array position+1)
# Values mirror row 2 ("1FDEU15H&K" / SYMBOL_2000_CA_SELECT_ENTRY_DATE case),
# varying only the BI/PD/UM/MP symbol columns (AE-AH), the entry date (AI)
# and the VALID flag (AJ) per renewal row.
$row3vals = @("1FDEU15H&K", "SYMBOL_2000_CA_SELECT_ENTRY_DATE", 2005, "CA_SE", "CA_SE", "Gt", "MDX ADVANCE", 53080, "WAG", "TEST", "SUV", "TEST", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", "K", 41, 41, "S", "Y", "N", "N", "N", "N", 20000101, "N", "Y", "N")
$row4vals = @("1FDEU15H&K", "SYMBOL_2000_CA_SELECT_ENTRY_DATE", 2005, "CA_SE", "CA_SE", "Gt", "MDX ADVANCE", 53080, "WAG", "TEST", "SUV", "TEST", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", "K", 41, 41, "S", "Y", "C", "C", "C", "C", 20150101, "Y", "Y", "N")
$row5vals = @("1FDEU15H&K", "SYMBOL_2000_CA_SELECT_ENTRY_DATE", 2005, "CA_SE", "CA_SE", "Gt", "MDX ADVANCE", 53080, "WAG", "TEST", "SUV", "TEST", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", "K", 41, 41, "S", "Y", "N", "N", "N", "N", 20180101, "N", "Y", "N")

$targetRows = @(3, 4, 5)
$allVals = @($row3vals, $row4vals, $row5vals)

# Shared-string insertion order matters for byte-identical output: the new
# "C" string (row 4's BI/PD/UM/MP symbol) must land in the table before the
# new "SYMBOL_2000_CA_SELECT_ENTRY_DATE" string (col B), so seed it first.
$ws.Cells.Item(4, 31).Value = "C"

for ($ri = 0; $ri -lt $targetRows.Length; $ri++) {
    $r = $targetRows[$ri]
    $vals = $allVals[$ri]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Re-apply the same formatting pattern used on row 2: col B keeps the
# "Good"-ish custom font style, most data columns keep the left-aligned
# style, and a handful (A, D, E, F, J, L, AI, AJ, AK, AL) stay default/unstyled.
$boldStyleCols = @(2)
$leftAlignCols = @(3, 7, 8, 9, 11, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34)
$noStyleCols = @(1, 4, 5, 6, 10, 12, 35, 36, 37, 38)

foreach ($r in $targetRows) {
    foreach ($c in $boldStyleCols) {
        $ws.Cells.Item(2, $c).Copy() | Out-Null
        $ws.Cells.Item($r, $c).PasteSpecial(-4122) | Out-Null
    }
    foreach ($c in $leftAlignCols) {
        $ws.Cells.Item(2, $c).Copy() | Out-Null
        $ws.Cells.Item($r, $c).PasteSpecial(-4122) | Out-Null
    }
    foreach ($c in $noStyleCols) {
        $ws.Cells.Item($r, $c).Style = "Normal"
    }
}

$excel.CutCopyMode = 0

# Column B is now wider text ("SYMBOL_2000_CA_SELECT_ENTRY_DATE") - autofit it.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Reset the view: selection on B5, scrolled back to show column A.
$ws.Range("A1").Select() | Out-Null
$ws.Range("B5").Select() | Out-Null
